# Auto-generated Excel COM-interop script
# Applies market-price data refresh to the Leve profit tables across all 8 job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the scheduled-runner update.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 930.5599999999999
$ws.Range("I41").Value = 1441.4546
$ws.Range("J41").Value = 529.1429000000001
$ws.Range("K41").Value = 1441.4546
$ws.Range("L41").Value = 529.1429000000001
$ws.Range("M41").Value = -1001.4546
$ws.Range("N41").Value = -1409.1429
# Row 125 (Leve Item ID 36228)
$ws.Range("H125").Value = 3052.2104
$ws.Range("I125").Value = 2573.7144
$ws.Range("K125").Value = 23163.4296
$ws.Range("M125").Value = -20703.4296
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 2138.8262
$ws.Range("I137").Value = 1552.9412
$ws.Range("K137").Value = 4658.8236
$ws.Range("M137").Value = -2108.8236

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 1488.5758
$ws.Range("I61").Value = 1411.138
$ws.Range("K61").Value = 1411.138
$ws.Range("M61").Value = -1199.138
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 3751.2
$ws.Range("I74").Value = 2287.4285
$ws.Range("J74").Value = 7166.6665
$ws.Range("K74").Value = 2287.4285
$ws.Range("L74").Value = 7166.6665
$ws.Range("M74").Value = -1413.4285
$ws.Range("N74").Value = -8914.666499999999
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 3751.2
$ws.Range("I77").Value = 2287.4285
$ws.Range("J77").Value = 7166.6665
$ws.Range("K77").Value = 11437.1425
$ws.Range("L77").Value = 35833.3325
$ws.Range("M77").Value = -7069.1425
$ws.Range("N77").Value = -44569.3325
# Row 88 (Leve Item ID 12530)
$ws.Range("H88").Value = 2711.2222
$ws.Range("I88").Value = 2475.25
$ws.Range("K88").Value = 2475.25
$ws.Range("M88").Value = -2069.25
# Row 91 (Leve Item ID 12530)
$ws.Range("H91").Value = 2711.2222
$ws.Range("I91").Value = 2475.25
$ws.Range("K91").Value = 2475.25
$ws.Range("M91").Value = -1071.25
# Row 125 (Leve Item ID 34251)
$ws.Range("H125").Value = 33850
$ws.Range("J125").Value = 33850
$ws.Range("L125").Value = 33850
$ws.Range("N125").Value = -43690
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 1428.1052
$ws.Range("I132").Value = 892.0769
$ws.Range("K132").Value = 2676.2307
$ws.Range("M132").Value = -146.2307000000001
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 1488.5758
$ws.Range("I136").Value = 1411.138
$ws.Range("K136").Value = 4233.414
$ws.Range("M136").Value = -1683.414

$ws = $wb.Worksheets.Item("BSM")
# Row 92 (Leve Item ID 18033)
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 1952.6552
$ws.Range("I134").Value = 1971.8077
$ws.Range("J134").Value = 1786.6666
$ws.Range("K134").Value = 5915.4231
$ws.Range("L134").Value = 5359.9998
$ws.Range("M134").Value = -3380.4231
$ws.Range("N134").Value = -10429.9998

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (Leve Item ID 27691)
$ws.Range("H16").Value = 2480.1428
$ws.Range("I16").Value = 2727.6667
$ws.Range("J16").Value = 995
$ws.Range("K16").Value = 2727.6667
$ws.Range("L16").Value = 995
$ws.Range("M16").Value = -2440.6667
$ws.Range("N16").Value = -1569
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 35705.08
$ws.Range("I31").Value = 25673.975
$ws.Range("J31").Value = 47169.2
$ws.Range("K31").Value = 25673.975
$ws.Range("L31").Value = 47169.2
$ws.Range("M31").Value = -25378.975
$ws.Range("N31").Value = -47759.2
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 35705.08
$ws.Range("I34").Value = 25673.975
$ws.Range("J34").Value = 47169.2
$ws.Range("K34").Value = 25673.975
$ws.Range("L34").Value = 47169.2
$ws.Range("M34").Value = -25471.975
$ws.Range("N34").Value = -47573.2
# Row 53 (Leve Item ID 25632)
$ws.Range("H53").Value = 40000
$ws.Range("J53").Value = 40000
$ws.Range("L53").Value = 40000
$ws.Range("N53").Value = -41214
# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 5082.372
$ws.Range("I58").Value = 1304.6086
$ws.Range("J58").Value = 9426.799999999999
$ws.Range("K58").Value = 1304.6086
$ws.Range("L58").Value = 9426.799999999999
$ws.Range("M58").Value = -1101.6086
$ws.Range("N58").Value = -9832.799999999999
# Row 86 (Leve Item ID 12584)
$ws.Range("H86").Value = 2185.524
$ws.Range("I86").Value = 1691.5385
$ws.Range("J86").Value = 2988.25
$ws.Range("K86").Value = 1691.5385
$ws.Range("L86").Value = 2988.25
$ws.Range("M86").Value = -568.5385000000001
$ws.Range("N86").Value = -5234.25
# Row 89 (Leve Item ID 12584)
$ws.Range("H89").Value = 2185.524
$ws.Range("I89").Value = 1691.5385
$ws.Range("J89").Value = 2988.25
$ws.Range("K89").Value = 8457.692500000001
$ws.Range("L89").Value = 14941.25
$ws.Range("M89").Value = -2841.692500000001
$ws.Range("N89").Value = -26173.25
# Row 96 (Leve Item ID 18193)
$ws.Range("H96").Value = 19155.75
$ws.Range("J96").Value = 19155.75
$ws.Range("L96").Value = 19155.75
$ws.Range("N96").Value = -24647.75
# Row 113 (Leve Item ID 27691)
$ws.Range("H113").Value = 2480.1428
$ws.Range("I113").Value = 2727.6667
$ws.Range("J113").Value = 995
$ws.Range("K113").Value = 2727.6667
$ws.Range("L113").Value = 995
$ws.Range("M113").Value = -557.6667000000002
$ws.Range("N113").Value = -5335
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2104.0125
$ws.Range("I132").Value = 1862.8448
$ws.Range("K132").Value = 5588.5344
$ws.Range("M132").Value = -3058.5344
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 960.931
$ws.Range("I134").Value = 946.3158
$ws.Range("K134").Value = 2838.9474
$ws.Range("M134").Value = -303.9474
# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 5082.372
$ws.Range("I136").Value = 1304.6086
$ws.Range("J136").Value = 9426.799999999999
$ws.Range("K136").Value = 3913.8258
$ws.Range("L136").Value = 28280.4
$ws.Range("M136").Value = -1363.8258
$ws.Range("N136").Value = -33380.39999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 1368.7
$ws.Range("I5").Value = 577.55554
$ws.Range("J5").Value = 3011.8462
$ws.Range("K5").Value = 1732.66662
$ws.Range("L5").Value = 9035.5386
$ws.Range("M5").Value = -1620.66662
$ws.Range("N5").Value = -9259.5386
# Row 34 (Leve Item ID 4749)
$ws.Range("H34").Value = 407.27274
$ws.Range("J34").Value = 651.6667
$ws.Range("L34").Value = 1955.0001
$ws.Range("N34").Value = -2123.0001
# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 773047.6
$ws.Range("I107").Value = 884.8333
$ws.Range("J107").Value = 1159129
$ws.Range("K107").Value = 2654.4999
$ws.Range("L107").Value = 3477387
$ws.Range("M107").Value = -734.4998999999998
$ws.Range("N107").Value = -3481227
# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 471.65216
$ws.Range("I113").Value = 463.66666
$ws.Range("J113").Value = 476.7857
$ws.Range("K113").Value = 1390.99998
$ws.Range("L113").Value = 1430.3571
$ws.Range("M113").Value = 779.0000199999999
$ws.Range("N113").Value = -5770.3571
# Row 122 (Leve Item ID 36078)
$ws.Range("H122").Value = 7219.1333
$ws.Range("I122").Value = 440.33334
$ws.Range("K122").Value = 3963.00006
$ws.Range("M122").Value = -1513.00006
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 840.0599999999999
$ws.Range("I131").Value = 610
$ws.Range("J131").Value = 854.7447
$ws.Range("K131").Value = 1830
$ws.Range("L131").Value = 2564.2341
$ws.Range("M131").Value = 3210
$ws.Range("N131").Value = -12644.2341
# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 1368.7
$ws.Range("I135").Value = 577.55554
$ws.Range("J135").Value = 3011.8462
$ws.Range("K135").Value = 5197.99986
$ws.Range("L135").Value = 27106.6158
$ws.Range("M135").Value = -2662.99986
$ws.Range("N135").Value = -32176.6158
# Row 140 (Leve Item ID 44097)
$ws.Range("H140").Value = 4379
$ws.Range("I140").Value = 4990.857
$ws.Range("K140").Value = 14972.571
$ws.Range("M140").Value = -9792.571

$ws = $wb.Worksheets.Item("GSM")
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 2794.6956
$ws.Range("I132").Value = 2286.2354
$ws.Range("K132").Value = 6858.706200000001
$ws.Range("M132").Value = -4328.706200000001
# Row 136 (Leve Item ID 42218)
$ws.Range("H136").Value = 14219.714
$ws.Range("J136").Value = 14219.714
$ws.Range("L136").Value = 42659.142
$ws.Range("N136").Value = -47759.142

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 2737.3157
$ws.Range("I7").Value = 1729.9
$ws.Range("J7").Value = 3856.6667
$ws.Range("K7").Value = 1729.9
$ws.Range("L7").Value = 3856.6667
$ws.Range("M7").Value = -1617.9
$ws.Range("N7").Value = -4080.6667
# Row 41 (Leve Item ID 3611)
$ws.Range("H41").Value = 6592.8887
$ws.Range("I41").Value = 3999
$ws.Range("J41").Value = 6917.125
$ws.Range("K41").Value = 3999
$ws.Range("L41").Value = 6917.125
$ws.Range("M41").Value = -3561
$ws.Range("N41").Value = -7793.125
# Row 42 (Leve Item ID 4333)
$ws.Range("H42").Value = 10137
$ws.Range("I42").Value = 8400
$ws.Range("J42").Value = 10426.5
$ws.Range("K42").Value = 8400
$ws.Range("L42").Value = 10426.5
$ws.Range("M42").Value = -7837
$ws.Range("N42").Value = -11552.5
# Row 49 (Leve Item ID 4333)
$ws.Range("H49").Value = 10137
$ws.Range("I49").Value = 8400
$ws.Range("J49").Value = 10426.5
$ws.Range("K49").Value = 8400
$ws.Range("L49").Value = 10426.5
$ws.Range("M49").Value = -8253
$ws.Range("N49").Value = -10720.5
# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 5480.6
$ws.Range("I68").Value = 2200
$ws.Range("J68").Value = 7667.6665
$ws.Range("K68").Value = 2200
$ws.Range("L68").Value = 7667.6665
$ws.Range("M68").Value = -1451
$ws.Range("N68").Value = -9165.666499999999
# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 5480.6
$ws.Range("I71").Value = 2200
$ws.Range("J71").Value = 7667.6665
$ws.Range("K71").Value = 11000
$ws.Range("L71").Value = 38338.3325
$ws.Range("M71").Value = -7256
$ws.Range("N71").Value = -45826.3325
# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 1465.3334
$ws.Range("I100").Value = 1558.4
$ws.Range("K100").Value = 1558.4
$ws.Range("M100").Value = -1017.4
# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 2737.3157
$ws.Range("I126").Value = 1729.9
$ws.Range("J126").Value = 3856.6667
$ws.Range("K126").Value = 5189.700000000001
$ws.Range("L126").Value = 11570.0001
$ws.Range("M126").Value = -2719.700000000001
$ws.Range("N126").Value = -16510.0001

$ws = $wb.Worksheets.Item("WVR")
# Row 92 (Leve Item ID 18088)
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0
# Row 135 (Leve Item ID 42043)
$ws.Range("H135").Value = 47599.6
$ws.Range("J135").Value = 47599.6
$ws.Range("L135").Value = 47599.6
$ws.Range("N135").Value = -57739.6
